$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Formula = '="67.398.14"'
$ws.Range('E2').Formula = '="  +4.89%  "'
$ws.Range('D3').Formula = '="3.259.92"'
$ws.Range('E3').Formula = '="  +2.74%  "'
$ws.Range('E4').Formula = '="  +0.02%  "'
$ws.Range('D5').Formula = '="579.06"'
$ws.Range('E5').Formula = '="  +2.53%  "'
$ws.Range('D6').Formula = '="179.65"'
$ws.Range('E6').Formula = '="  +5.78%  "'
$ws.Range('E7').Formula = '="  -0.02%  "'
$ws.Range('D8').Formula = '="0.601"'
$ws.Range('E8').Formula = '="  -0.82%  "'
$ws.Range('D9').Formula = '="3.258.36"'
$ws.Range('E9').Formula = '="  +2.79%  "'
$ws.Range('E10').Formula = '="  +4.17%  "'
$ws.Range('D11').Formula = '="6.75"'
$ws.Range('E11').Formula = '="  +1.63%  "'
$ws.Range('E12').Formula = '="  +4.55%  "'
$ws.Range('D13').Formula = '="3.823.83"'
$ws.Range('E13').Formula = '="  +2.65%  "'
$ws.Range('E14').Formula = '="  +0.47%  "'
$ws.Range('D15').Formula = '="28.22"'
$ws.Range('E15').Formula = '="  +3.15%  "'
$ws.Range('D16').Formula = '="67.383.60"'
$ws.Range('E16').Formula = '="  +4.88%  "'
$ws.Range('D17').Formula = '="0.0000169"'
$ws.Range('E17').Formula = '="  +3.05%  "'
$ws.Range('D18').Formula = '="3.258.19"'
$ws.Range('E18').Formula = '="  +2.65%  "'
$ws.Range('D19').Formula = '="5.89"'
$ws.Range('E19').Formula = '="  +2.54%  "'
$ws.Range('D20').Formula = '="13.42"'
$ws.Range('E20').Formula = '="  +3.35%  "'
$ws.Range('D21').Formula = '="377.96"'
$ws.Range('E21').Formula = '="  +7.12%  "'
$ws.Range('D22').Formula = '="7.66"'
$ws.Range('E22').Formula = '="  +6.66%  "'
$ws.Range('E23').Formula = '="  +0.10%  "'
$ws.Range('D24').Formula = '="71.41"'
$ws.Range('E24').Formula = '="  +3.51%  "'
$ws.Range('E25').Formula = '="  +1.84%  "'
$ws.Range('D26').Formula = '="3.397.66"'
$ws.Range('E26').Formula = '="  +2.53%  "'
$ws.Range('E27').Formula = '="  -1.04%  "'
$ws.Range('D28').Formula = '="9.91"'
$ws.Range('E28').Formula = '="  +3.77%  "'
$ws.Range('E29').Formula = '="  +1.82%  "'
$ws.Range('E30').Formula = '="  +0.06%  "'
$ws.Range('E31').Formula = '="  +4.24%  "'
$ws.Range('D32').Formula = '="5.65"'
$ws.Range('E32').Formula = '="  +1.54%  "'
$ws.Range('D33').Formula = '="22.64"'
$ws.Range('E33').Formula = '="  +2.72%  "'
$ws.Range('E34').Formula = '="  +0.10%  "'
$ws.Range('E35').Formula = '="  +5.93%  "'
$ws.Range('D36').Formula = '="6.84"'
$ws.Range('E36').Formula = '="  +3.23%  "'
$ws.Range('D37').Formula = '="164.05"'
$ws.Range('E37').Formula = '="  +6.42%  "'
$ws.Range('E38').Formula = '="  +4.26%  "'
$ws.Range('D39').Formula = '="0.861"'
$ws.Range('E39').Formula = '="  +5.33%  "'
$ws.Range('E40').Formula = '="  +9.65%  "'
$ws.Range('D41').Formula = '="27.04"'
$ws.Range('E41').Formula = '="  +4.44%  "'
$ws.Range('B42').Formula = '="RenderToken"'
$ws.Range('C42').Formula = '="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"'
$ws.Range('D42').Formula = '="6.74"'
$ws.Range('E42').Formula = '="  +11.45%  "'
$ws.Range('B43').Formula = '="dogwifhat"'
$ws.Range('C43').Formula = '="https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"'
$ws.Range('D43').Formula = '="2.64"'
$ws.Range('E43').Formula = '="  +3.64%  "'
$ws.Range('D44').Formula = '="2.775.61"'
$ws.Range('E44').Formula = '="  +5.71%  "'
$ws.Range('D45').Formula = '="4.40"'
$ws.Range('E45').Formula = '="  +5.25%  "'
$ws.Range('D46').Formula = '="25.86"'
$ws.Range('D47').Formula = '="347.24"'
$ws.Range('E47').Formula = '="  +8.09%  "'
$ws.Range('D48').Formula = '="40.53"'
$ws.Range('E48').Formula = '="  +2.55%  "'
$ws.Range('D49').Formula = '="0.0676"'
$ws.Range('E49').Formula = '="  +3.13%  "'
$ws.Range('E50').Formula = '="  +3.66%  "'
$ws.Range('E51').Formula = '="  +1.26%  "'

$usedRange = $ws.UsedRange
$usedRange.Copy()
$usedRange.PasteSpecial(-4163)
$excel.CutCopyMode = 0
